# Auto-generated edit script: updates Siren Profits sheet values
# per the commit diff (scheduled-runner price/profit recompute).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3670.7727
$ws.Range("I18").Value = 3807.476
$ws.Range("K18").Value = 3807.476
$ws.Range("M18").Value = -3523.476
$ws.Range("H51").Value = 33236.668
$ws.Range("J51").Value = 36735.938
$ws.Range("L51").Value = 36735.938
$ws.Range("N51").Value = -37703.938
$ws.Range("H62").Value = 47625676
$ws.Range("I62").Value = 62507536
$ws.Range("J62").Value = 3717.2
$ws.Range("K62").Value = 62507536
$ws.Range("L62").Value = 3717.2
$ws.Range("M62").Value = -62506912
$ws.Range("N62").Value = -4965.2
$ws.Range("H65").Value = 47625676
$ws.Range("I65").Value = 62507536
$ws.Range("J65").Value = 3717.2
$ws.Range("K65").Value = 312537680
$ws.Range("L65").Value = 18586
$ws.Range("M65").Value = -312534560
$ws.Range("N65").Value = -24826
$ws.Range("H70").Value = 5128.9375
$ws.Range("J70").Value = 4251.769
$ws.Range("L70").Value = 12755.307
$ws.Range("N70").Value = -13295.307
$ws.Range("H73").Value = 5128.9375
$ws.Range("J73").Value = 4251.769
$ws.Range("L73").Value = 12755.307
$ws.Range("N73").Value = -14627.307
$ws.Range("H98").Value = 18687.457
$ws.Range("I98").Value = 22761.32
$ws.Range("J98").Value = 8502.799999999999
$ws.Range("K98").Value = 22761.32
$ws.Range("L98").Value = 8502.799999999999
$ws.Range("M98").Value = -21263.32
$ws.Range("N98").Value = -11498.8
$ws.Range("H107").Value = 9583
$ws.Range("I107").Value = 11456.444
$ws.Range("J107").Value = 1152.5
$ws.Range("K107").Value = 11456.444
$ws.Range("L107").Value = 1152.5
$ws.Range("M107").Value = -9536.444
$ws.Range("N107").Value = -4992.5
$ws.Range("H122").Value = 18687.457
$ws.Range("I122").Value = 22761.32
$ws.Range("J122").Value = 8502.799999999999
$ws.Range("K122").Value = 68283.95999999999
$ws.Range("L122").Value = 25508.4
$ws.Range("M122").Value = -65833.95999999999
$ws.Range("N122").Value = -30408.4
$ws.Range("H137").Value = 290745.1
$ws.Range("I137").Value = 413194.1
$ws.Range("K137").Value = 1239582.3
$ws.Range("M137").Value = -1237032.3
$ws.Range("H138").Value = 5269.75
$ws.Range("I138").Value = 2296.9
$ws.Range("J138").Value = 5977.5713
$ws.Range("K138").Value = 6890.700000000001
$ws.Range("L138").Value = 17932.7139
$ws.Range("M138").Value = -1750.700000000001
$ws.Range("N138").Value = -28212.7139

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3127.3408
$ws.Range("I32").Value = 3165.325
$ws.Range("K32").Value = 3165.325
$ws.Range("M32").Value = -2878.325
$ws.Range("H61").Value = 17707.824
$ws.Range("I61").Value = 20502.785
$ws.Range("K61").Value = 20502.785
$ws.Range("M61").Value = -20290.785
$ws.Range("H74").Value = 1524.7838
$ws.Range("I74").Value = 1412.2646
$ws.Range("K74").Value = 1412.2646
$ws.Range("M74").Value = -538.2646
$ws.Range("H77").Value = 1524.7838
$ws.Range("I77").Value = 1412.2646
$ws.Range("K77").Value = 7061.323
$ws.Range("M77").Value = -2693.323
$ws.Range("H132").Value = 2121.4307
$ws.Range("I132").Value = 1749.7258
$ws.Range("K132").Value = 5249.1774
$ws.Range("M132").Value = -2719.1774
$ws.Range("H136").Value = 17707.824
$ws.Range("I136").Value = 20502.785
$ws.Range("K136").Value = 61508.355
$ws.Range("M136").Value = -58958.355

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3389.8333
$ws.Range("I20").Value = 2195.9285
$ws.Range("K20").Value = 2195.9285
$ws.Range("M20").Value = -1948.9285
$ws.Range("H134").Value = 5039.8223
$ws.Range("I134").Value = 5041.6743
$ws.Range("K134").Value = 15125.0229
$ws.Range("M134").Value = -12590.0229

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 11112470
$ws.Range("I7").Value = 2787.25
$ws.Range("J7").Value = 20000216
$ws.Range("K7").Value = 2787.25
$ws.Range("L7").Value = 20000216
$ws.Range("M7").Value = -2674.25
$ws.Range("N7").Value = -20000442
$ws.Range("H18").Value = 125000
$ws.Range("J18").Value = 125000
$ws.Range("L18").Value = 125000
$ws.Range("N18").Value = -125460
$ws.Range("H22").Value = 12820687
$ws.Range("I22").Value = 199.5
$ws.Range("J22").Value = 38461660
$ws.Range("K22").Value = 199.5
$ws.Range("L22").Value = 38461660
$ws.Range("M22").Value = 150.5
$ws.Range("N22").Value = -38462360
$ws.Range("H31").Value = 3201.3901
$ws.Range("I31").Value = 1551
$ws.Range("J31").Value = 3601.4849
$ws.Range("K31").Value = 1551
$ws.Range("L31").Value = 3601.4849
$ws.Range("M31").Value = -1256
$ws.Range("N31").Value = -4191.484899999999
$ws.Range("H34").Value = 3201.3901
$ws.Range("I34").Value = 1551
$ws.Range("J34").Value = 3601.4849
$ws.Range("K34").Value = 1551
$ws.Range("L34").Value = 3601.4849
$ws.Range("M34").Value = -1349
$ws.Range("N34").Value = -4005.4849
$ws.Range("H58").Value = 7396.7744
$ws.Range("I58").Value = 11295.8
$ws.Range("J58").Value = 3741.4375
$ws.Range("K58").Value = 11295.8
$ws.Range("L58").Value = 3741.4375
$ws.Range("M58").Value = -11092.8
$ws.Range("N58").Value = -4147.4375
$ws.Range("H99").Value = 316466.25
$ws.Range("I99").Value = 717785.7
$ws.Range("J99").Value = 4328.8887
$ws.Range("K99").Value = 717785.7
$ws.Range("L99").Value = 4328.8887
$ws.Range("M99").Value = -716287.7
$ws.Range("N99").Value = -7324.8887
$ws.Range("H117").Value = 50712
$ws.Range("J117").Value = 50712
$ws.Range("L117").Value = 50712
$ws.Range("N117").Value = -59890
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()  # was -64056
$ws.Range("H121").Value = 21864
$ws.Range("I121").Value = 20296
$ws.Range("K121").Value = 20296
$ws.Range("M121").Value = -18986
$ws.Range("H126").Value = 316466.25
$ws.Range("I126").Value = 717785.7
$ws.Range("J126").Value = 4328.8887
$ws.Range("K126").Value = 2153357.1
$ws.Range("L126").Value = 12986.6661
$ws.Range("M126").Value = -2150887.1
$ws.Range("N126").Value = -17926.6661
$ws.Range("H136").Value = 7396.7744
$ws.Range("I136").Value = 11295.8
$ws.Range("J136").Value = 3741.4375
$ws.Range("K136").Value = 33887.39999999999
$ws.Range("L136").Value = 11224.3125
$ws.Range("M136").Value = -31337.39999999999
$ws.Range("N136").Value = -16324.3125

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2554.0256
$ws.Range("J68").Value = 2655.0625
$ws.Range("L68").Value = 7965.1875
$ws.Range("N68").Value = -9587.1875
$ws.Range("H71").Value = 2554.0256
$ws.Range("J71").Value = 2655.0625
$ws.Range("L71").Value = 23895.5625
$ws.Range("N71").Value = -32007.5625
$ws.Range("H121").Value = 1177713.9
$ws.Range("I121").Value = 232.90909
$ws.Range("K121").Value = 698.72727
$ws.Range("M121").Value = 611.27273
$ws.Range("H131").Value = 4656.8
$ws.Range("J131").Value = 1907.85
$ws.Range("L131").Value = 5723.549999999999
$ws.Range("N131").Value = -15803.55
$ws.Range("H137").Value = 2967.926
$ws.Range("I137").Value = 2557.8
$ws.Range("J137").Value = 8094.5
$ws.Range("K137").Value = 7673.400000000001
$ws.Range("L137").Value = 24283.5
$ws.Range("M137").Value = -2573.400000000001
$ws.Range("N137").Value = -34483.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2727.7144
$ws.Range("I113").Value = 1297.5
$ws.Range("K113").Value = 1297.5
$ws.Range("H132").Value = 3396.4375
$ws.Range("I132").Value = 3109.1316
$ws.Range("J132").Value = 4488.2
$ws.Range("K132").Value = 9327.3948
$ws.Range("L132").Value = 13464.6
$ws.Range("M132").Value = -6797.3948
$ws.Range("N132").Value = -18524.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 40057.082
$ws.Range("I7").Value = 46593.5
$ws.Range("K7").Value = 46593.5
$ws.Range("M7").Value = -46481.5
$ws.Range("H22").Value = 16520.941
$ws.Range("I22").Value = 18941.1
$ws.Range("K22").Value = 18941.1
$ws.Range("M22").Value = -18646.1
$ws.Range("H27").Value = 16520.941
$ws.Range("I27").Value = 18941.1
$ws.Range("K27").Value = 18941.1
$ws.Range("M27").Value = -18834.1
$ws.Range("H40").Value = 45031.89
$ws.Range("I40").Value = 60432
$ws.Range("K40").Value = 60432
$ws.Range("M40").Value = -60296
$ws.Range("H126").Value = 40057.082
$ws.Range("I126").Value = 46593.5
$ws.Range("K126").Value = 139780.5
$ws.Range("M126").Value = -137310.5
$ws.Range("H132").Value = 10337.972
$ws.Range("I132").Value = 15479.632
$ws.Range("K132").Value = 46438.896
$ws.Range("M132").Value = -43908.896

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 15243.875
$ws.Range("J41").Value = 14408.5
$ws.Range("L41").Value = 14408.5
$ws.Range("N41").Value = -15188.5
$ws.Range("H51").Value = 57498.5
$ws.Range("I51").Value = 56666.332
$ws.Range("K51").Value = 56666.332
$ws.Range("M51").Value = -56156.332
$ws.Range("H107").Value = 2152.2083
$ws.Range("I107").Value = 1797.7778
$ws.Range("J107").Value = 3215.5
$ws.Range("K107").Value = 5393.3334
$ws.Range("L107").Value = 9646.5
$ws.Range("M107").Value = -3473.3334
$ws.Range("N107").Value = -13486.5
$ws.Range("H113").Value = 2573.919
$ws.Range("I113").Value = 905.4074000000001
$ws.Range("J113").Value = 7078.9
$ws.Range("K113").Value = 2716.2222
$ws.Range("L113").Value = 21236.7
$ws.Range("M113").Value = -546.2222000000002
$ws.Range("N113").Value = -25576.7
